$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New value for the previously-empty C2 cell
$ws.Range("C2").Value = 3.09825757489699

# Updated values for existing cells in columns C and E (rows 2-19)
$values = @{
    2  = @{ E = 7.865470614547343 }
    3  = @{ C = -1.791203563722299; E = -6.760862998203643 }
    4  = @{ C = 0.5799958470386946; E = 5.643342995751777 }
    5  = @{ C = 3.181454202131073;  E = 4.38978860149748 }
    6  = @{ C = 0.5930547804883668; E = -1.194610791900008 }
    7  = @{ C = -0.3951783438669754; E = 0.03694906323863378 }
    8  = @{ C = 3.292216014290039;  E = 7.617133650412211 }
    9  = @{ C = 1.670328650030184;  E = 2.037906845818616 }
    10 = @{ C = 2.562791874943371;  E = 3.265947405805814 }
    11 = @{ C = 1.526411006965533;  E = 0.6601843988560674 }
    12 = @{ C = 1.63465618619294;   E = 1.551857746372698 }
    13 = @{ C = 1.35261353265177;   E = 0.8024032016000104 }
    14 = @{ C = -2.082763426755907; E = -5.866344937500023 }
    15 = @{ C = -0.1380317107957718; E = 7.749494937649115 }
    16 = @{ C = 3.848999231984762;  E = 2.866003071127765 }
    17 = @{ C = -0.3745803349312071; E = 0.645722451525943 }
    18 = @{ C = -1.091476630333243; E = -0.4907904687545206 }
    19 = @{ C = 1.626992717807862;  E = 0.6270138473519316 }
}

foreach ($row in $values.Keys) {
    $cols = $values[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
